$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.624.64'
$ws.Range('E2').Value = '  -1.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.677.41'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.38'
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3930'
$ws.Range('E7').Value = '  -2.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3959'
$ws.Range('E8').Value = '  -2.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.001'
$ws.Range('E9').Value = '  -0.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.407'
$ws.Range('E10').Value = '  -3.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '50.97'
$ws.Range('E11').Value = '  -5.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08671'
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '25.39'
$ws.Range('E13').Value = '  -1.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.347'
$ws.Range('E14').Value = '  -1.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001323'
$ws.Range('E15').Value = '  -1.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.732'
$ws.Range('E16').Value = '  -3.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.662.78'
$ws.Range('E17').Value = '  -7.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.94'
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.54'
$ws.Range('E19').Value = '  +2.79%  '
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.07028'
$ws.Range('E20').Value = '  -2.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.090'
$ws.Range('E21').Value = '  -2.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9992'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.99'
$ws.Range('E23').Value = '  -4.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.634.49'
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.368'
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.780'
$ws.Range('E26').Value = '  -3.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.17'
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.887'
$ws.Range('E28').Value = '  -10.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '160.79'
$ws.Range('E29').Value = '  -1.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '146.86'
$ws.Range('E30').Value = '  +2.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.334'
$ws.Range('E31').Value = '  +2.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.515'
$ws.Range('E32').Value = '  +10.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.852.50'
$ws.Range('E33').Value = '  -6.60%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.03115'
$ws.Range('E34').Value = '  -2.15%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08322'
$ws.Range('E35').Value = '  -4.94%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.974'
$ws.Range('E36').Value = '  -5.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2815'
$ws.Range('E37').Value = '  -0.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9908'
$ws.Range('E38').Value = '  -3.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.09536'
$ws.Range('E39').Value = '  +1.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.516'
$ws.Range('E40').Value = '  +3.14%  '
$ws.Range('E41').Value = '  -5.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7930'
$ws.Range('E42').Value = '  -6.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.63'
$ws.Range('E43').Value = '  -2.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.61'
$ws.Range('E44').Value = '  -6.54%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7144'
$ws.Range('E45').Value = '  -4.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.565'
$ws.Range('E46').Value = '  -5.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.168'
$ws.Range('E47').Value = '  -1.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.08661'
$ws.Range('E48').Value = '  +3.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.9994'
$ws.Range('E49').Value = '  -0.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.334'
$ws.Range('E50').Value = '  -5.24%  '
$ws.Range('E51').Value = '  -2.60%  '
